$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.353.45'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').Value = '2.651.44'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.72'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('E6').Value = '  +2.75%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.143'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +3.75%  '
$ws.Range('E10').Value = '  -1.21%  '
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('E12').Value = '  +0.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.13'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.55%  '
$ws.Range('E14').Value = '  +1.33%  '
$ws.Range('D15').Value = '3.131.67'
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = '68.313.20'
$ws.Range('E16').Value = '  +0.16%  '
$ws.Range('D17').Value = '2.661.25'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '364.39'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.40%  '
$ws.Range('E20').Value = '  -0.80%  '
$ws.Range('E21').Value = '  +3.78%  '
$ws.Range('E22').Value = '  -0.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.07'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '75.24'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.23%  '
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('E26').Value = '  -2.83%  '
$ws.Range('E27').Value = '  +0.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000105'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.54%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '558.89'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.03'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('E32').Value = '  -0.17%  '
$ws.Range('E33').Value = '  +0.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.129'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.82%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('E36').Value = '  +1.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.84'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '159.70'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.45%  '
$ws.Range('E39').Value = '  +1.06%  '
$ws.Range('E40').Value = '  -2.14%  '
$ws.Range('E41').Value = '  -0.37%  '
$ws.Range('E42').Value = '  +3.33%  '
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '158.49'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.00%  '
$ws.Range('E46').Value = '  +0.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '22.18'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.48%  '
$ws.Range('E48').Value = '  -0.83%  '
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.615'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.567'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.12%  '
